$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.146.47"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "1.802.30"
$ws.Range("E3").Value = "  +3.42%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'335.87"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").Value = "'0.9991"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").Value = "'0.4666"
$ws.Range("E7").Value = "  +24.16%  "
$ws.Range("D8").Value = "'0.3711"
$ws.Range("E8").Value = "  +10.73%  "
$ws.Range("D9").Value = "'45.34"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  +6.55%  "
$ws.Range("D11").Value = "'1.153"
$ws.Range("E11").Value = "  +3.54%  "
$ws.Range("D12").Value = "'22.60"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").Value = "'6.381"
$ws.Range("E14").Value = "  +3.64%  "
$ws.Range("D15").Value = "'7.410"
$ws.Range("E15").Value = "  +4.32%  "
$ws.Range("D16").Value = "1.799.24"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").Value = "'0.00001095"
$ws.Range("E17").Value = "  +3.62%  "
$ws.Range("D18").Value = "'0.06726"
$ws.Range("D19").Value = "'82.63"
$ws.Range("E19").Value = "  +4.07%  "
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").Value = "'17.44"
$ws.Range("E21").Value = "  +3.68%  "
$ws.Range("D22").Value = "'6.426"
$ws.Range("E22").Value = "  +3.07%  "
$ws.Range("D23").Value = "28.156.28"
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("D24").Value = "'11.92"
$ws.Range("E24").Value = "  +2.40%  "
$ws.Range("D25").Value = "'2.416"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").Value = "'20.86"
$ws.Range("D27").Value = "'2.397"
$ws.Range("E27").Value = "  +3.66%  "
$ws.Range("D28").Value = "'152.47"
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("D29").Value = "2.005.87"
$ws.Range("E29").Value = "  +3.02%  "
$ws.Range("D30").Value = "'134.58"
$ws.Range("E30").Value = "  +2.35%  "
$ws.Range("E31").Value = "  +2.01%  "
$ws.Range("D32").Value = "'4.041"
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("D33").Value = "'0.09640"
$ws.Range("E33").Value = "  +10.66%  "
$ws.Range("D34").Value = "'5.920"
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("E35").Value = "  +6.72%  "
$ws.Range("D36").Value = "'12.24"
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("D37").Value = "'0.02374"
$ws.Range("E37").Value = "  +2.60%  "
$ws.Range("D38").Value = "'0.06396"
$ws.Range("E38").Value = "  +3.22%  "
$ws.Range("D39").Value = "'0.6729"
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("D40").Value = "'5.272"
$ws.Range("E40").Value = "  +2.53%  "
$ws.Range("D41").Value = "'1.522"
$ws.Range("E41").Value = "  +6.01%  "
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("D43").Value = "'8.125"
$ws.Range("E43").Value = "  +2.11%  "
$ws.Range("D44").Value = "'14.20"
$ws.Range("E44").Value = "  +2.99%  "
$ws.Range("D45").Value = "'0.9988"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("D46").Value = "'0.6184"
$ws.Range("E46").Value = "  +2.44%  "
$ws.Range("D47").Value = "'3.840"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("D48").Value = "'129.92"
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("E49").Value = "  +2.60%  "
$ws.Range("D50").Value = "'1.187"
$ws.Range("E50").Value = "  +0.75%  "
$ws.Range("D51").Value = "'0.07145"
